$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.5023523333333334
$ws.Range("H2").Value = 1.507057
$ws.Range("I2").Value = 0.06515888850144765
$ws.Range("J2").Value = 0.06515888850144765
$ws.Range("M2").Value = 0.08241233333333334
$ws.Range("Q2").Value = 0.04140002794544445
$ws.Range("R2").Value = 0.3726002515090001
$ws.Range("S2").Value = 0.003080064316741594
$ws.Range("T2").Value = 0.003080064316741594
$ws.Range("G3").Value = 0.5023523333333334
$ws.Range("H3").Value = 1.507057
$ws.Range("I3").Value = 0.06515888850144765
$ws.Range("J3").Value = 0.06515888850144765
$ws.Range("Q3").Value = 0.8344192821226668
$ws.Range("R3").Value = 7.509773539104001
$ws.Range("S3").Value = 0.06207882418470605
$ws.Range("T3").Value = 0.06207882418470605
$ws.Range("I4").Value = 0.07275905893716338
$ws.Range("J4").Value = 0.07275905893716339
$ws.Range("M4").Value = 0.08241233333333334
$ws.Range("Q4").Value = 0.04622895114633333
$ws.Range("S4").Value = 0.003439324799824917
$ws.Range("T4").Value = 0.003439324799824918
$ws.Range("I5").Value = 0.07275905893716338
$ws.Range("J5").Value = 0.07275905893716339
$ws.Range("S5").Value = 0.06931973413733847
$ws.Range("T5").Value = 0.06931973413733847
$ws.Range("G6").Value = 2.845667666666667
$ws.Range("H6").Value = 8.537003
$ws.Range("I6").Value = 0.3691045704399529
$ws.Range("J6").Value = 0.3691045704399529
$ws.Range("M6").Value = 0.08241233333333334
$ws.Range("Q6").Value = 0.2345181123012222
$ws.Range("R6").Value = 2.110663010711
$ws.Range("S6").Value = 0.01744759376202489
$ws.Range("T6").Value = 0.01744759376202489
$ws.Range("G7").Value = 2.845667666666667
$ws.Range("H7").Value = 8.537003
$ws.Range("I7").Value = 0.3691045704399529
$ws.Range("J7").Value = 0.3691045704399529
$ws.Range("Q7").Value = 4.726722290357334
$ws.Range("R7").Value = 42.540500613216
$ws.Range("S7").Value = 0.351656976677928
$ws.Range("T7").Value = 0.3516569766779279
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.11543
$ws.Range("H8").Value = 0.34629
$ws.Range("I8").Value = 0.01497214206175765
$ws.Range("J8").Value = 0.01497214206175765
$ws.Range("M8").Value = 0.08241233333333334
$ws.Range("Q8").Value = 0.009512855636666666
$ws.Range("R8").Value = 0.08561570073000001
$ws.Range("S8").Value = 0.000707733995624881
$ws.Range("T8").Value = 0.0007077339956248811
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.11543
$ws.Range("H9").Value = 0.34629
$ws.Range("I9").Value = 0.01497214206175765
$ws.Range("J9").Value = 0.01497214206175765
$ws.Range("Q9").Value = 0.19173200032
$ws.Range("R9").Value = 1.72558800288
$ws.Range("S9").Value = 0.01426440806613277
$ws.Range("T9").Value = 0.01426440806613277
$ws.Range("G10").Value = 3.685254666666667
$ws.Range("H10").Value = 11.055764
$ws.Range("I10").Value = 0.4780053400596784
$ws.Range("J10").Value = 0.4780053400596784
$ws.Range("M10").Value = 0.08241233333333334
$ws.Range("Q10").Value = 0.3037104360075556
$ws.Range("R10").Value = 2.733393924068
$ws.Range("S10").Value = 0.02259533925439868
$ws.Range("T10").Value = 0.02259533925439868
$ws.Range("G11").Value = 3.685254666666667
$ws.Range("H11").Value = 11.055764
$ws.Range("I11").Value = 0.4780053400596784
$ws.Range("J11").Value = 0.4780053400596784
$ws.Range("Q11").Value = 6.121296447445334
$ws.Range("R11").Value = 55.091668027008
$ws.Range("S11").Value = 0.4554100008052797
$ws.Range("T11").Value = 0.4554100008052797
